# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the "conversion" note text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$text = $cellA1.Value()
$text = $text.Replace("1000 Bs = 3.28 = 12520.36 pesos", "1000 Bs = 3.28 = 12563.13 pesos")
$text = $text.Replace("12520.36 pesos = 3.26 = 954.83 Bs", "12563.13 pesos = 3.28 = 960.09 Bs")
$cellA1.Value = $text

# --- Sheet "tasas": refresh the rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("O10").Value = 3833.05
$wsTasas.Range("N12").Value = 3834
